# Staging.Role.xlsx - "moved staging files StagingTemplates directory"
#
# The underlying column headers were re-ordered (RoleID/Code/Name/Description
# -> Code/Description/Name/RoleID) and the per-column "best fit" width
# overrides that were tailored to the old header text are dropped so the
# sheet falls back to the default column width for columns B:D (column A
# keeps its explicit width). The workbook/sheet also picked up fresh
# window-size / VBA codeName bookkeeping from being re-saved in its new
# location - we set what the object model exposes for that too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header re-order -------------------------------------------------
# Row 2 holds the column headers. Swap them into their new positions:
#   A2: RoleID -> Code
#   B2: Code   -> Description
#   C2: Name   -> Name        (unchanged)
#   D2: Description -> RoleID
#
# Insert three blank, unformatted columns at B:D (this pushes the old
# B:D - with their bespoke bestFit widths - out to E:G without touching
# column A), then delete the now-displaced old columns outright so their
# width overrides go away instead of just being reassigned. What's left
# at B:D are plain, default-width columns ready for the new header text.
$ws.Columns("B:D").Insert()
$ws.Columns("E:G").Delete()

$ws.Range("A2").Value = "Code"
$ws.Range("B2").Value = "Description"
$ws.Range("C2").Value = "Name"
$ws.Range("D2").Value = "RoleID"

# --- Cosmetic workbook/sheet bookkeeping ------------------------------
# Best-effort: mirror the saved window size and sheet VBA codeName from
# the commit. (No-ops if the host doesn't expose these for scripting.)
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 28800
    $win.Height = 12585
} catch {}

try {
    $ws.CodeName = "Sheet48"
} catch {}
